$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows for 15, 16, 17 (aggiornamento 15, 16, 17 marzo)
$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 6
$ws.Range("C227").Value = 56
$ws.Range("D227").Value = 162.9517546412151

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 7
$ws.Range("C228").Value = 47
$ws.Range("D228").Value = 136.7630797881627

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 8
$ws.Range("C229").Value = 48
$ws.Range("D229").Value = 139.672932549613

# Carry over the date-column formatting (bold, bordered, centered, date numfmt)
# from the last existing row so the new cells match the established style.
$ws.Range("A226").Copy()
$ws.Range("A227:A229").PasteSpecial(-4122)
